# Update computed result values on several report sheets.
$wb = $excel.ActiveWorkbook

# --- pcroprep (sheet1) ---
$ws1 = $wb.Worksheets.Item("pcroprep")
$ws1.Range("D35").Value = [Convert]::ToDouble("2.3845454553868214E-14")
$ws1.Range("F35").Value = [Convert]::ToDouble("8.671074383224803E-14")
$ws1.Range("G35").Value = [Convert]::ToDouble("-241.39999999999992")
$ws1.Range("D39").Value = [Convert]::ToDouble("736.03566632590366")
$ws1.Range("F39").Value = [Convert]::ToDouble("225.63973114196455")
$ws1.Range("G39").Value = [Convert]::ToDouble("-580.96026885803553")

# --- pdietrep (sheet4) ---
$ws4 = $wb.Worksheets.Item("pdietrep")
$ws4.Range("E6").Value = [Convert]::ToDouble("954.61828080308487")
$ws4.Range("F6").Value = [Convert]::ToDouble("-1203.167233402598")
$ws4.Range("G6").Value = [Convert]::ToDouble("44.240647391429725")
$ws4.Range("E7").Value = [Convert]::ToDouble("33.530307718074909")
$ws4.Range("F7").Value = [Convert]::ToDouble("-36.59788315827565")
$ws4.Range("G7").Value = [Convert]::ToDouble("47.812879954646583")
$ws4.Range("E8").Value = [Convert]::ToDouble("12.186140051228314")
$ws4.Range("F8").Value = [Convert]::ToDouble("-52.547425374942158")
$ws4.Range("G8").Value = [Convert]::ToDouble("18.825071616249495")
$ws4.Range("E9").Value = [Convert]::ToDouble("171.50702489674933")
$ws4.Range("F9").Value = [Convert]::ToDouble("-152.16080223410302")
$ws4.Range("G9").Value = [Convert]::ToDouble("52.988592167800633")

# --- pradar (sheet5) ---
$ws5 = $wb.Worksheets.Item("pradar")
$ws5.Range("D15").Value = [Convert]::ToDouble("8.671074383224803E-14")
$ws5.Range("E15").Value = [Convert]::ToDouble("3.5919943592480544E-14")
$ws5.Range("F15").Value = [Convert]::ToDouble("-241.39999999999992")

# --- plandrep (sheet6) ---
$ws6 = $wb.Worksheets.Item("plandrep")
$ws6.Range("S11").Value = [Convert]::ToDouble("1.7638685125293154E-13")

# --- plaborrep (sheet7) ---
$ws7 = $wb.Worksheets.Item("plaborrep")
$ws7.Range("R3").Value = [Convert]::ToDouble("1.4375528377113891E-16")
$ws7.Range("AF3").Value = [Convert]::ToDouble("0.82435132156501723")

# --- pfertrep (sheet8) ---
$ws8 = $wb.Worksheets.Item("pfertrep")
$ws8.Range("S5").Value = [Convert]::ToDouble("4.8506384094556175E-11")
$ws8.Range("Z5").Value = [Convert]::ToDouble("190785.88409599999")
$ws8.Range("S6").Value = [Convert]::ToDouble("2.9103830456733704E-11")
$ws8.Range("Z6").Value = [Convert]::ToDouble("218764.48264999996")
$ws8.Range("S7").Value = [Convert]::ToDouble("4.6389741879520996E-11")
$ws8.Range("Z7").Value = [Convert]::ToDouble("231384.044096")
